# Inserts a new price-report row for Terminal Hortofrutícola Agro Chillán - Tomate
# at row 213 (pushing the previous rows 213-303 down to 214-304), and populates
# the new row with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 213; everything below (213-303) shifts down to 214-304.
$ws.Rows("213:213").Insert()

# Populate the newly inserted row 213 with the new observation.
$ws.Range("A213").Value = 7
$ws.Range("B213").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C213").Value = "Ñuble"
$ws.Range("D213").Value = 44489
$ws.Range("E213").Value = 16
$ws.Range("F213").Value = 100112020
$ws.Range("G213").Value = "Tomate"
$ws.Range("H213").Value = "Larga vida"
$ws.Range("I213").Value = "Primera"
$ws.Range("J213").Value = 400
$ws.Range("K213").Value = 7500
$ws.Range("L213").Value = 8000
$ws.Range("M213").Value = 7750
$ws.Range("N213").Value = '$/caja 10 kilos'
$ws.Range("O213").Value = "Región de Arica y Parinacota"
$ws.Range("P213").Value = 775
$ws.Range("Q213").Value = 10
$ws.Range("R213").Value = "Hortaliza"
